# Add a new "2021" data column (column Y) to the indicator table, mirroring
# the formatting of the existing "2020" column (column X), and update the
# sheet's active selection to the cell that was selected when the workbook
# was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds the year header values; row 5-16 hold one data series per row.
# For each row, copy the 2020 cell's formatting (number format, font,
# border, alignment) into the new 2021 cell in column Y, then overwrite it
# with the 2021 value so the copied 2020 value doesn't stick around.
$newColumnValues = [ordered]@{
    4  = 2021
    5  = 46.69
    6  = 52.52
    7  = 43.22
    8  = 51.31
    9  = 41.31
    10 = 52.43
    11 = 49.27
    12 = 31.68
    13 = 35.59
    14 = 55.28
    15 = 61.02
    16 = 48.72
}

foreach ($row in $newColumnValues.Keys) {
    $srcCell = $ws.Range("X$row")
    $dstCell = $ws.Range("Y$row")
    $srcCell.Copy($dstCell)
    $dstCell.Value = $newColumnValues[$row]
}

# Reflect the saved cursor/selection position recorded in the workbook view.
$ws.Range("AA15").Select()
